# Add "RIGHT" option fields for "hungry" and "sleepy", mirroring the
# existing Pet/Noise (SINGLE) rows, on the Info sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# --- Hungry state/wants (row 9-10) ---
$ws.Range("A9").Value = "Hungry"
$ws.Range("B9").Value = "SPREADSHEETFORM:RIGHT:hungry:state"

$ws.Range("A10").Value = "Wants"
$ws.Range("B10").Value = "SPREADSHEETFORM:RIGHT:hungry:wants"

# --- Sleepy state/wants (row 13-14) ---
$ws.Range("A13").Value = "Sleepy"
$ws.Range("B13").Value = "SPREADSHEETFORM:RIGHT:sleepy:state"

$ws.Range("A14").Value = "Wants"
$ws.Range("B14").Value = "SPREADSHEETFORM:RIGHT:sleepy:wants"

# Copy the existing label / value formatting (bold-with-border for column A,
# italic-with-border for column B) from the first pair of rows so the new
# rows look consistent with the existing ones.
$ws.Range("A5").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New "scratch" columns C/D/E next to the new rows, boxed with a thin border
# (same border as the existing bordered cells) but left blank / default font.
$ws.Range("C9:E10").Borders.LineStyle = 1
$ws.Range("C13:E14").Borders.LineStyle = 1

# Give the 3 new columns explicit widths (~22.62 / ~22.51 / ~23.2 chars).
$ws.Columns.Item(3).ColumnWidth = 21.833333
$ws.Columns.Item(4).ColumnWidth = 21.666667
$ws.Columns.Item(5).ColumnWidth = 22.333333

# Two blank spacer rows below the new blocks (rows 16-17), same row height
# as the rest of the sheet.
$ws.Rows.Item(16).RowHeight = 15
$ws.Rows.Item(17).RowHeight = 15

# Restore the selection like the target workbook (cursor parked below the
# newly added content).
$ws.Range("B19").Select()
